# The deck's Date placeholders (on every slide layout, every slide
# master, and the notes master) cache the rendered value of an
# auto-updating datetime field. The presentation was reopened/saved a
# few days later, so every cached "10/23/2024" needs to become
# "10/27/2024".

$p = $ppt.ActivePresentation

$oldDate = "10/23/2024"
$newDate = "10/27/2024"

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

function Update-DateInShapes($shapes) {
    for ($si = 1; $si -le $shapes.Count; $si++) {
        Update-DateShape $shapes.Item($si)
    }
}

for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $design = $p.Designs.Item($di)
    $master = $design.SlideMaster

    # The master's own Date placeholder.
    Update-DateInShapes $master.Shapes

    if ($di -eq 1) {
        # Layouts that belong to the first (in-use) master can be
        # reached directly and reliably through the master.
        $layouts = $master.CustomLayouts
        for ($li = 1; $li -le $layouts.Count; $li++) {
            Update-DateInShapes $layouts.Item($li).Shapes
        }
    } else {
        # Layouts belonging to any additional master aren't used by
        # any slide yet, and addressing them straight off
        # Master.CustomLayouts.Item(n) isn't reliable here, so borrow a
        # temporary slide, swap its CustomLayout across each of this
        # master's layouts in turn (which resolves them correctly),
        # edit, then discard the helper slide.
        $tempSlide = $p.Slides.Add($p.Slides.Count + 1, 1)
        $layouts = $master.CustomLayouts
        for ($li = 1; $li -le $layouts.Count; $li++) {
            $tempSlide.CustomLayout = $layouts.Item($li)
            Update-DateInShapes $tempSlide.CustomLayout.Shapes
        }
        $tempSlide.Delete()
    }
}

# The notes master's Date placeholder.
Update-DateInShapes $p.NotesMaster.Shapes
